# Pandoc reference.docx style-sheet update ("Add Figure (from Pandoc 3)")
#
# Adds the "Abstract Title" paragraph style (and re-tunes "Abstract"'s
# spacing), adds the "Footnote Block Text" paragraph style, and fills in
# the previously-empty ImportTok / BuiltInTok syntax-highlighting
# character styles.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New "Abstract Title" style, inserted ahead of "Abstract" in the
#    document flow (Next -> Abstract).
# ---------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.Alignment = 1

$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 0x8A5A34  # wdColor BGR encoding of RGB 345A8A
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10

# ---------------------------------------------------------------------
# 2. "Abstract" keeps its look, but SpaceBefore tightens from 15pt
#    (300 twips) to 5pt (100 twips); SpaceAfter stays 15pt (300 twips).
# ---------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5
$abstract.ParagraphFormat.SpaceAfter = 15

# ---------------------------------------------------------------------
# 3. New "Footnote Block Text" style, based on Footnote Text.
# ---------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0

# ---------------------------------------------------------------------
# 4. Fill in syntax-highlighting character styles that used to be empty.
# ---------------------------------------------------------------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Bold = $true
$importTok.Font.Color = 0x008000

$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 0x008000
